# Revised semantics of properties resource; Removed propName attribute of Job
# This script removes the "PROPERTIES propName" clause from the PUT JOB and
# UPDATE JOB command syntax (column A), removes the corresponding
# jobs/{jobName}/propfile URI + PUT/DELETE method rows from the UPDATE JOB
# table cell (columns C/D on row 22), shrinks the now-shorter wrapped rows
# to their new heights, and restores the active selection to A2 on the
# Commands sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# Row 19: PUT JOB jobName SCRIPT scriptName ...
# Remove the "[NO PROPERTIES | PROPERTIES propName]" line.
$ws.Range("A19").Value = "PUT JOB jobName SCRIPT scriptName" + "`n" + `
    "[NO ARGUMENTS | ARGUMENTS argName1 argValue1 [, " + [char]0x2026 + "]]" + "`n" + `
    "[NO SCHEDULE[S] | SCHEDULE[S] schedName1 [, " + [char]0x2026 + "]]" + "`n" + `
    "[ENABLED {ON | OFF}]"

# Row 22: UPDATE JOB jobName [SCRIPT scriptName] ...
# Remove the "[NO PROPERTIES |" / "PROPERTIES propName]" lines.
$ws.Range("A22").Value = "UPDATE JOB jobName [SCRIPT scriptName]" + "`n" + `
    "[NO ARGUMENTS |" + "`n" + `
    "ARGUMENTS argName1 argValue1 [, " + [char]0x2026 + "]]" + "`n" + `
    "[NO SCHEDULE[S] |" + "`n" + `
    "SCHEDULE[S] schedName1 [, " + [char]0x2026 + "]]" + "`n" + `
    "[ENABLED {ON | OFF}]"

# Row 22 Method column: drop the DELETE/PUT pair for the removed propfile URI.
$ws.Range("C22").Value = "PUT`nDELETE`nPUT`nPUT"

# Row 22 URI Template column: drop the two jobs/{jobName}/propfile lines.
$ws.Range("D22").Value = "jobs/{jobName}/script`njobs/{jobName}/arguments`njobs/{jobName}/arguments`njobs/{jobName}/schedules`njobs/{jobName}/schedules`njobs/{jobName}/enabled"

# The cell text is now shorter (fewer wrapped lines), so the rows shrink.
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 90

# Restore the active cell/selection to A2 on the Commands sheet.
$ws.Activate()
$ws.Range("A2").Select()
